$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: number of clusters goes from 2 to 1 (row now matches the GatesS/GatesT
# "single cluster / no meaningful split" rows below it).
$ws.Range("B2").Value = 1

# C2/D2 (avg within-cluster variance / silhouette score) become empty text
# cells, exactly like C3/D3 below. A bare ".Value = ''' would make Excel
# clear the cell entirely, so we write the text-prefix marker "'" (forces an
# empty *text* cell instead of a blank one) and then restore the default
# "Normal" style so no stray quote-prefix/number-format sticks to the cell.
$ws.Range("C2").Value = "'"
$ws.Range("D2").Value = "'"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Style = "Normal"

# E2/F2: most significant feature / p-value are no longer meaningful, so they
# report "N/A" like the other rows.
$ws.Range("E2").Value = "N/A"
$ws.Range("F2").Value = "N/A"
